# feat: add 2022-Q3 data
#
# - Insert a new worksheet "2022-Q3" right after "总计" (becomes the 2nd tab,
#   pushing "2021-Q1" and "2020-Q4" one position to the right).
# - Populate "2022-Q3" with the same fund-holding table shape used by the
#   other quarterly sheets.
# - Add a summary row for 2022-Q3 into "总计" (as the new row 2), shifting the
#   existing 2021-Q1 / 2020-Q4 summary rows down by one row.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item("总计")
$q1_2021 = $wb.Worksheets.Item("2021-Q1")

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q3" worksheet right after "总计".
# ---------------------------------------------------------------------------
$q3_2022 = $wb.Worksheets.Add($null, $total)
$q3_2022.Name = "2022-Q3"

# Header row
$q3_2022.Range("B1").Value = "基金代码"
$q3_2022.Range("C1").Value = "基金名称"
$q3_2022.Range("D1").Value = "基金规模"
$q3_2022.Range("E1").Value = "股票总仓位"
$q3_2022.Range("F1").Value = "仓位占比"
$q3_2022.Range("G1").Value = "持有市值(亿元)"
$q3_2022.Range("H1").Value = "仓位排名"

# Row 2
$q3_2022.Range("A2").Value = 0
$q3_2022.Range("B2").Value = "'005585"
$q3_2022.Range("C2").Value = "银河文体娱乐主题灵活配置混合A"
$q3_2022.Range("D2").Value = "'3.01"
$q3_2022.Range("E2").Value = "'90.28"
$q3_2022.Range("F2").Value = "'3.76"
$q3_2022.Range("G2").Value = "'0.1132"
$q3_2022.Range("H2").Value = 10

# Row 3
$q3_2022.Range("A3").Value = 1
$q3_2022.Range("B3").Value = "'015667"
$q3_2022.Range("C3").Value = "银河文体娱乐主题灵活配置混合C"
$q3_2022.Range("D3").Value = "'0.41"
$q3_2022.Range("E3").Value = "'90.28"
$q3_2022.Range("F3").Value = "'3.76"
$q3_2022.Range("G3").Value = "'0.0154"
$q3_2022.Range("H3").Value = 10

# Match the header / id-column styling used by the sibling quarterly sheets
# (bold font + border + centered alignment) by copying the formatting over.
$q1_2021.Range("B1:H1").Copy()
$q3_2022.Range("B1:H1").PasteSpecial(-4122)
$q1_2021.Range("A2:A3").Copy()
$q3_2022.Range("A2:A3").PasteSpecial(-4122)

$q3_2022.Range("A1").Select()

# ---------------------------------------------------------------------------
# 2. Shift the "总计" summary rows down and insert the 2022-Q3 summary.
# ---------------------------------------------------------------------------
# Push 2020-Q4 (currently row 3) down to row 4 - copy values then formats so
# the destination cells end up identical to the source (incl. A-column style).
$total.Range("A3:D3").Copy()
$total.Range("A4:D4").PasteSpecial(-4163)
$total.Range("A3:D3").Copy()
$total.Range("A4:D4").PasteSpecial(-4122)
$total.Range("A4").Value = 2

# Push 2021-Q1 (currently row 2) down to row 3.
$total.Range("A2:D2").Copy()
$total.Range("A3:D3").PasteSpecial(-4163)
$total.Range("A2:D2").Copy()
$total.Range("A3:D3").PasteSpecial(-4122)
$total.Range("A3").Value = 1

# Write the new 2022-Q3 summary into row 2 (row already carries the right
# formatting, just overwrite the values).
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.13

$total.Range("A1").Select()
